$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Table style change (the "PLENARY" recap table on slide 16):
#    {7006D89D-E2AB-4669-9D1A-4C46F1BE7EEB} -> {E65A3A18-3622-46E5-A44E-AEC817A4C587}
#    Walk every slide/shape defensively and restyle whichever table(s)
#    are still using the old style id.
# ------------------------------------------------------------------
$oldStyleId = "{7006D89D-E2AB-4669-9D1A-4C46F1BE7EEB}"
$newStyleId = "{E65A3A18-3622-46E5-A44E-AEC817A4C587}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            if ($shp.Table.Style -eq $oldStyleId) {
                $shp.Table.ApplyStyle($newStyleId)
            }
        }
    }
}

# ------------------------------------------------------------------
# 2) Theme swap: the deck's theme ("Integral") and the colours that
#    used to live in the secondary theme part ("Office Theme") trade
#    places. The slide master's active theme colour scheme is driven
#    through the object model here, updating every themed colour slot
#    (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) from the
#    "Integral" palette to the "Office" palette.
# ------------------------------------------------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB = 0            # dk1     000000
$tcs.Item(2).RGB = 16777215     # lt1     FFFFFF
$tcs.Item(3).RGB = 6968388      # dk2     44546A
$tcs.Item(4).RGB = 15132391     # lt2     E7E6E6
$tcs.Item(5).RGB = 13998939     # accent1 5B9BD5
$tcs.Item(6).RGB = 3243501      # accent2 ED7D31
$tcs.Item(7).RGB = 10855845     # accent3 A5A5A5
$tcs.Item(8).RGB = 49407        # accent4 FFC000
$tcs.Item(9).RGB = 12874308     # accent5 4472C4
$tcs.Item(10).RGB = 4697456     # accent6 70AD47
$tcs.Item(11).RGB = 12673797    # hlink   0563C1
$tcs.Item(12).RGB = 7491477     # folHlink 954F72
